$d = $word.ActiveDocument

$d.Content.Find.Execute("How does the application of the computer science help doctors in reading the psychological behaviour of the Humans and help them by measuring the personality traits to get rid from the unusual behaviour?", $true, $false, $false, $false, $false, $true, 1, $false, "How does the application of computer science help doctors in reading the psychological behaviour of Humans and help them by measuring personality traits to get rid from the unusual behaviour?", 2) | Out-Null
$d.Content.Find.Execute("How the application of computer science helping the people to track their health statistics?", $true, $false, $false, $false, $false, $true, 1, $false, "How the application of computer science helps people track their health statistics?", 2) | Out-Null
$d.Content.Find.Execute("At the moment the advancement in the computer field is exponentially high, in my opinion since the computers became easily accessible and easy to handle I feel if we make use of this technology for the benefit of the humans in medical field where it is needed to collect the big amount of data for the research and to provide treatment based on the results, as it is with medical field the errors in findings may result giving wrong treatment to people and this may impact on their life which is quite sensitive to deal with, to make sure there shouldn’t be any errors in the collected data which is in real difficult to handle manually, if we collect data digitally or if we teach the machines to get the reviews based on the actions and how he the person reacts at times and keep a track of the activities which is with so many people may help to get correct information and helps to get the results without intervein of humans, which will help the doctors to have an idea for the root cause and do research on the collected data which is complex data structure, in this case if we collaborate the computer science to the medical field which will ease doctors to get results accurately and handle the cases correctly, for instance talking about psychology to research on people with the events like what making people to act in a different way and under what circumstances made big impact to their personality trait, for this they have to observe with different people and need to collect the data like their activities and the situations happened in their lives, for this they need to collect from a lot of people which is quite hard to deal by humans as needed to deal with big amount of data, in this kind of scenario if we apply the computer science technology embedded to the medical field that may help doctors to understand the case so clearly without errors.", $true, $false, $false, $false, $false, $true, 1, $false, "The advancement in the field of computer is exponential, in my opinion since computers became easily accessible and easy to handle if we make use of this technology for the benefit of humans in the medical field where it is needed to collect the big amount of data for research and to provide treatment based on the results, as it is with medical field the errors in findings may result giving wrong treatment to people and this may impact on their life which is quite sensitive to deal with, to make sure there shouldn’t be any errors in the collected data which is in real difficult to handle manually, if we collect data digitally or we teach the machines to get the reviews based on the actions and how he the person reacts at times and keep a track of the activities of people may help to get correct information and helps to get the results without intervein of humans, which will help the doctors to have an idea for the root cause and do research on the collected data which is complex data structure, in this case collaboration the computer science and the medical field which will ease doctors workload to get results accurately and handle the cases correctly, for instance talking about psychology to research on people with the events like what making people to act in a different way and under what circumstances made big impact to their personality trait, for this they have to observe with different people and need to collect the data like their activities and the situations happened in their lives, for this they need to collect from a lot of sources which is quite hard to deal by humans as needed to deal with big amount of data, in this kind of scenario if we apply the computer science technology embedded to the medicine such statistical data may improve accuracy to treat. ", 2) | Out-Null
$d.Content.Find.Execute("As per the health statistics 1 in 4 people having the mental issues and also with mental problem like depression, anxiety. This demonstrates that there is necessity of improvement in that area as the expected numbers are high, for this we need sophisticated technology which we can train machines to deal with generating results to repeated issues, we can make it possible by connecting the computer science technology which is INTERNET OF THINGS to the medical field the analytics can be tracked by connecting the device to the central server, for example we can track the hours spent by user in browsing the kind of data and the history demonstrates the actual behaviour of the user by connecting the data transfer to the central server.", $true, $false, $false, $false, $false, $true, 1, $false, "As per the health statistics 1 in 4 people are diagnosed with mental problem like depression, anxiety. This demonstrates that there is necessity of improvement in that area as the expected numbers are high, for this we need sophisticated technology which we can train machines to deal with generating results to repeated issues, we can make it possible by connecting the computer science technology which is INTERNET OF THINGS to the medical field the analytics can be tracked by connecting the device to the central server, for example we can track the hours spent by user in browsing the kind of data and the history demonstrates the actual behaviour of the user by connecting the data transfer to the central server.", 2) | Out-Null
$d.Content.Find.Execute("To handle this big data, we need technologies to handle the data like accessing, storing, processing and retrieving the desired data from the server, we definitely need Big data technology and well efficient servers to store and access the data quickly.", $true, $false, $false, $false, $false, $true, 1, $false, "To handle this big data, we need technologies for accessing, storing, processing and retrieving the desired data from the server, we definitely need Big data technology and well efficient servers to store and access the data quickly.", 2) | Out-Null
$d.Content.Find.Execute("By having the values, the collaboration with the computer science and psychology how the body reacts to any external thing with the signal to the brain, heart rate, which lead to more innovative inventions like the car stops if the driver is sleepy, and also learning software recognise if a student drops attention.", $true, $false, $false, $false, $false, $true, 1, $false, "By having values of collaboration between computer science and psychology how the body reacts to any external thing with the signal to the brain, heart rate, which lead to more innovative inventions like the car stops if the driver is sleepy, and also learning software recognise if a student drops attention.", 2) | Out-Null
$d.Content.Find.Execute("The collect the data it is important to interact with the people of different personalities, tacking their actions against the corresponding incidents, it is imperative to have experiments on people with questionnaires, reviews and sometimes interviews.", $true, $false, $false, $false, $false, $true, 1, $false, "Collecting data is important for interacting with people of different personalities, tracking their actions against the corresponding incidents, it is imperative to have experiments on people with questionnaires, reviews and sometimes interviews.", 2) | Out-Null
$d.Content.Find.Execute("Without applying the computer science technology to the medical field it might be hard to do certain treatments to people, as it is important to have sensors based to see what is happening inside the body also with the values.", $true, $false, $false, $false, $false, $true, 1, $false, "Without applying computer science technology to the medicine, it might be hard to treat to people, as it is important to have sensors based to see what is happening inside the body also with values.", 2) | Out-Null
$d.Content.Find.Execute("Making possibility to ease and to get accurate results to help the people to overcome sophisticated things.", $true, $false, $false, $false, $false, $true, 1, $false, "Making possibility to ease and to get accurate results to help the people to get the right services at times.", 2) | Out-Null
